# Fix the misspelling "Shakespear" -> "Shakespeare".
#
# B2 holds the text used as the "Text" label for the first block of rows
# (rows 2-4). It is referenced by the shared-string table, by the summary
# formula in B15 (=B2), and (via B15) by the cached category labels in the
# "conclusion" chart. Updating the cell's value is enough to fix the
# shared string table and the dependent formula; the embedded chart title
# that also spells out the word is fixed separately below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the misspelled source text.
$ws.Range("B2").Value = "Shakespeare"

# Recalculate so the B15 (=B2) summary formula and other dependents pick
# up the corrected text.
$excel.CalculateFull()

# The first line chart's title is static rich text ("CPU" + " TIme
# Shakespear") rather than a cell reference, so it has to be corrected
# directly as well.
$co = $ws.ChartObjects(1)
$chart = $co.Chart()
$chart.ChartTitle().Text = "CPU TIme Shakespeare"
